# Apply corrected IFRS financial figures to company_list sheet (rows 2-9)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 138537
$ws.Range("E2").Value = 13883
$ws.Range("F2").Value = 13883
$ws.Range("G2").Value = 13433
$ws.Range("H2").Value = 10320
$ws.Range("I2").Value = 10266
$ws.Range("J2").Value = 54
$ws.Range("K2").Value = 2197608
$ws.Range("L2").Value = 2041145
$ws.Range("M2").Value = 156463
$ws.Range("N2").Value = 155623
$ws.Range("O2").Value = 840
$ws.Range("P2").Value = 32559
$ws.Range("Q2").Value = -33705
$ws.Range("R2").Value = 12097
$ws.Range("S2").Value = 10598
$ws.Range("T2").Value = 0
$ws.Range("V2").Value = 1003896
$ws.Range("W2").Value = 10.02
$ws.Range("X2").Value = 7.45
$ws.Range("Y2").Value = 6.86
$ws.Range("Z2").Value = 0.48
$ws.Range("AA2").Value = 1304.56
$ws.Range("AB2").Value = 380.54
$ws.Range("AC2").Value = 1578
$ws.Range("AD2").Value = 8.93
$ws.Range("AE2").Value = 23898
$ws.Range("AF2").Value = 0.59
$ws.Range("AG2").Value = 430
$ws.Range("AH2").Value = 3.05
$ws.Range("AI2").Value = 27.28
$ws.Range("AJ2").Value = 553217245
$ws.Range("U2").ClearContents()

# Row 3
$ws.Range("D3").Value = 137327
$ws.Range("E3").Value = 14997
$ws.Range("F3").Value = 14997
$ws.Range("G3").Value = 14695
$ws.Range("H3").Value = 11506
$ws.Range("I3").Value = 11430
$ws.Range("J3").Value = 76
$ws.Range("K3").Value = 2398428
$ws.Range("L3").Value = 2225490
$ws.Range("M3").Value = 172938
$ws.Range("N3").Value = 172030
$ws.Range("O3").Value = 908
$ws.Range("P3").Value = 32719
$ws.Range("Q3").Value = -60070
$ws.Range("R3").Value = -10134
$ws.Range("S3").Value = 66536
$ws.Range("V3").Value = 1062824
$ws.Range("W3").Value = 10.92
$ws.Range("X3").Value = 8.38
$ws.Range("Y3").Value = 6.98
$ws.Range("Z3").Value = 0.5
$ws.Range("AA3").Value = 1286.87
$ws.Range("AB3").Value = 428.56
$ws.Range("AC3").Value = 1747
$ws.Range("AD3").Value = 7.07
$ws.Range("AE3").Value = 26289
$ws.Range("AF3").Value = 0.47
$ws.Range("AG3").Value = 450
$ws.Range("AH3").Value = 3.64
$ws.Range("AI3").Value = 25.76
$ws.Range("AJ3").Value = 556401958
$ws.Range("T3").ClearContents()
$ws.Range("U3").ClearContents()

# Row 4
$ws.Range("D4").Value = 152777
$ws.Range("E4").Value = 15326
$ws.Range("F4").Value = 15326
$ws.Range("G4").Value = 15172
$ws.Range("H4").Value = 11646
$ws.Range("I4").Value = 11575
$ws.Range("J4").Value = 71
$ws.Range("K4").Value = 2568514
$ws.Range("L4").Value = 2388040
$ws.Range("M4").Value = 180474
$ws.Range("N4").Value = 179500
$ws.Range("O4").Value = 974
$ws.Range("P4").Value = 32898
$ws.Range("Q4").Value = -40981
$ws.Range("R4").Value = -12191
$ws.Range("S4").Value = 54346
$ws.Range("V4").Value = 1147378
$ws.Range("W4").Value = 10.03
$ws.Range("X4").Value = 7.62
$ws.Range("Y4").Value = 6.59
$ws.Range("Z4").Value = 0.47
$ws.Range("AA4").Value = 1323.2
$ws.Range("AB4").Value = 448.59
$ws.Range("AC4").Value = 1762
$ws.Range("AD4").Value = 7.21
$ws.Range("AE4").Value = 27282
$ws.Range("AF4").Value = 0.47
$ws.Range("AG4").Value = 480
$ws.Range("AH4").Value = 3.78
$ws.Range("AI4").Value = 27.28
$ws.Range("AJ4").Value = 559978815
$ws.Range("T4").ClearContents()
$ws.Range("U4").ClearContents()

# Row 5
$ws.Range("D5").Value = 167840
$ws.Range("E5").Value = 20283
$ws.Range("F5").Value = 20283
$ws.Range("G5").Value = 19536
$ws.Range("H5").Value = 15085
$ws.Range("I5").Value = 15015
$ws.Range("J5").Value = 71
$ws.Range("K5").Value = 2740697
$ws.Range("L5").Value = 2542680
$ws.Range("M5").Value = 198017
$ws.Range("N5").Value = 196974
$ws.Range("O5").Value = 1044
$ws.Range("P5").Value = 32898
$ws.Range("Q5").Value = 31985
$ws.Range("R5").Value = -77726
$ws.Range("S5").Value = 44240
$ws.Range("V5").Value = 1190230
$ws.Range("W5").Value = 12.08
$ws.Range("X5").Value = 8.99
$ws.Range("Y5").Value = 7.98
$ws.Range("Z5").Value = 0.57
$ws.Range("AA5").Value = 1284.07
$ws.Range("AB5").Value = 501.92
$ws.Range("AC5").Value = 2282
$ws.Range("AD5").Value = 7.21
$ws.Range("AE5").Value = 29937
$ws.Range("AF5").Value = 0.55
$ws.Range("AG5").Value = 617
$ws.Range("AH5").Value = 3.75
$ws.Range("AI5").Value = 27.04
$ws.Range("AJ5").Value = 559978815
$ws.Range("T5").ClearContents()
$ws.Range("U5").ClearContents()

# Row 6
$ws.Range("D6").Value = 155279
$ws.Range("E6").Value = 23964
$ws.Range("F6").Value = 23964
$ws.Range("G6").Value = 23995
$ws.Range("H6").Value = 17643
$ws.Range("I6").Value = 17542
$ws.Range("K6").Value = 2895094
$ws.Range("L6").Value = 2684153
$ws.Range("M6").Value = 210942
$ws.Range("N6").Value = 209829
$ws.Range("P6").Value = 32898
$ws.Range("Q6").Value = -103969
$ws.Range("R6").Value = -40544
$ws.Range("S6").Value = 161970
$ws.Range("V6").Value = 1359548
$ws.Range("W6").Value = 15.43
$ws.Range("X6").Value = 11.36
$ws.Range("Y6").Value = 8.63
$ws.Range("Z6").Value = 0.63
$ws.Range("AA6").Value = 1272.46
$ws.Range("AB6").Value = 541.21
$ws.Range("AC6").Value = 2666
$ws.Range("AD6").Value = 5.27
$ws.Range("AE6").Value = 31891
$ws.Range("AF6").Value = 0.44
$ws.Range("AG6").Value = 690
$ws.Range("AH6").Value = 4.91
$ws.Range("AI6").Value = 23.38
$ws.Range("AJ6").Value = 559978815
$ws.Range("T6").ClearContents()
$ws.Range("U6").ClearContents()

# Row 7
$ws.Range("E7").Value = 25173
$ws.Range("G7").Value = 23014
$ws.Range("H7").Value = 16846
$ws.Range("I7").Value = 16812
$ws.Range("K7").Value = 3091896
$ws.Range("L7").Value = 2864021
$ws.Range("M7").Value = 227501
$ws.Range("N7").Value = 227518
$ws.Range("P7").Value = 33580
$ws.Range("Y7").Value = 7.69
$ws.Range("Z7").Value = 0.56
$ws.Range("AA7").Value = 1258.9
$ws.Range("AC7").Value = 2505
$ws.Range("AD7").Value = 4.29
$ws.Range("AE7").Value = 33700
$ws.Range("AF7").Value = 0.32
$ws.Range("AG7").Value = 688
$ws.Range("AH7").Value = 6.4
$ws.Range("AI7").Value = 23.6
$ws.Range("D7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()

# Row 8
$ws.Range("E8").Value = 25752
$ws.Range("G8").Value = 23000
$ws.Range("H8").Value = 16766
$ws.Range("I8").Value = 16763
$ws.Range("K8").Value = 3254580
$ws.Range("L8").Value = 3014906
$ws.Range("M8").Value = 239441
$ws.Range("N8").Value = 239062
$ws.Range("P8").Value = 33580
$ws.Range("Y8").Value = 7.18
$ws.Range("Z8").Value = 0.53
$ws.Range("AA8").Value = 1259.14
$ws.Range("AC8").Value = 2483
$ws.Range("AD8").Value = 4.33
$ws.Range("AE8").Value = 35410
$ws.Range("AF8").Value = 0.3
$ws.Range("AG8").Value = 697
$ws.Range("AH8").Value = 6.48
$ws.Range("AI8").Value = 23.99
$ws.Range("D8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()

# Row 9
$ws.Range("E9").Value = 29690
$ws.Range("G9").Value = 23776
$ws.Range("H9").Value = 17463
$ws.Range("I9").Value = 17345
$ws.Range("K9").Value = 3373105
$ws.Range("L9").Value = 3120570
$ws.Range("M9").Value = 252538
$ws.Range("N9").Value = 251232
$ws.Range("P9").Value = 33518
$ws.Range("Y9").Value = 7.07
$ws.Range("Z9").Value = 0.53
$ws.Range("AA9").Value = 1235.69
$ws.Range("AC9").Value = 2569
$ws.Range("AD9").Value = 4.18
$ws.Range("AE9").Value = 37212
$ws.Range("AF9").Value = 0.29
$ws.Range("AG9").Value = 736
$ws.Range("AH9").Value = 6.85
$ws.Range("AI9").Value = 24.5
$ws.Range("D9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
